# "add stats vs VSS" — round 8 results for Д3 tournament.
# Updates the standings table (rows 5-18) with the new games-played /
# wins / losses / score / points numbers, and appends the two new match
# days (2025-01-18 and 2025-01-19) with their game results at the bottom
# of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Standings table (rows 5-18): Игры(D) / Побед(E) / Поражений(F) /
#    Команда(C) / Мячи(G) / Очки(H) after round 8.
# ---------------------------------------------------------------------

$standings = @(
    @{ Row = 5;  Team = "ISsoft";                 Games = 8; Wins = 6; Losses = 2; Score = "515 - 438"; Points = 14 },
    @{ Row = 6;  Team = "Эра-Недвижимости плюс";   Games = 8; Wins = 6; Losses = 2; Score = "604 - 497"; Points = 14 },
    @{ Row = 7;  Team = "ОПЛАТИ";                  Games = 8; Wins = 6; Losses = 2; Score = "573 - 514"; Points = 14 },
    @{ Row = 8;  Team = "GOLDEN HILL";             Games = 8; Wins = 6; Losses = 2; Score = "564 - 522"; Points = 14 },
    @{ Row = 9;  Team = "Грушвиль";                Games = 8; Wins = 6; Losses = 2; Score = "630 - 493"; Points = 14 },
    @{ Row = 10; Team = "Mapogo males";            Games = 8; Wins = 6; Losses = 2; Score = "611 - 550"; Points = 14 },
    @{ Row = 11; Team = "БГУФК";                   Games = 8; Wins = 5; Losses = 3; Score = "569 - 456"; Points = 13 },
    @{ Row = 12; Team = "SIRIUS";                  Games = 8; Wins = 5; Losses = 3; Score = "542 - 434"; Points = 13 },
    @{ Row = 13; Team = "VSS";                     Games = 8; Wins = 3; Losses = 5; Score = "480 - 522"; Points = 11 },
    @{ Row = 14; Team = "Стрела";                  Games = 8; Wins = 2; Losses = 6; Score = "471 - 530"; Points = 10 },
    @{ Row = 15; Team = "NORD";                    Games = 8; Wins = 2; Losses = 6; Score = "407 - 638"; Points = 10 },
    @{ Row = 16; Team = "Eagles";                  Games = 8; Wins = 2; Losses = 6; Score = "468 - 494"; Points = 10 },
    @{ Row = 17; Team = "ЛФК";                     Games = 8; Wins = 1; Losses = 7; Score = "436 - 564"; Points = 9 },
    @{ Row = 18; Team = "Минск 7х";                Games = 8; Wins = 0; Losses = 8; Score = "361 - 579"; Points = 8 }
)

foreach ($s in $standings) {
    $r = $s.Row
    $ws.Range("C$r").Value = $s.Team
    $ws.Range("D$r").Value = $s.Games
    $ws.Range("E$r").Value = $s.Wins
    $ws.Range("F$r").Value = $s.Losses
    $ws.Range("G$r").Value = $s.Score
    $ws.Range("H$r").Value = $s.Points
}

# ---------------------------------------------------------------------
# 2. New rows 84-92: the two match days played after the last existing
#    entry (row 83), copying the formatting of the last existing date
#    block (rows 79-83) and filling in the new schedule text.
# ---------------------------------------------------------------------

# Row 84 — date header for 2025-01-18 (serial 45675), formatted like row 79.
$ws.Range("B79:H79").Copy()
$ws.Range("B84:H84").PasteSpecial(-4122)
$ws.Range("B84").Value = 45675

# Rows 85-87 — three game results for 2025-01-18, formatted like row 80
# (copy/paste one row at a time so every cell in the row gets the same
# style instead of the paste being tiled across a multi-row block).
$ws.Range("B80:H80").Copy()
$ws.Range("B85:H85").PasteSpecial(-4122)
$ws.Range("B80:H80").Copy()
$ws.Range("B86:H86").PasteSpecial(-4122)
$ws.Range("B80:H80").Copy()
$ws.Range("B87:H87").PasteSpecial(-4122)

$ws.Rows.Item(85).RowHeight = 19.95
$ws.Rows.Item(86).RowHeight = 19.95
$ws.Rows.Item(87).RowHeight = 19.95

$ws.Range("B85").Value = "Грушвиль - Минск 7х 110:32 (16:30, БНТУ)"
$ws.Range("B86").Value = "БГУФК - Стрела 72:44 (18:00, БНТУ)"
$ws.Range("B87").Value = "ОПЛАТИ - Eagles 84:77 (19:30, БНТУ)"

# Row 88 — date header for 2025-01-19 (serial 45676), formatted like row 79.
$ws.Range("B79:H79").Copy()
$ws.Range("B88:H88").PasteSpecial(-4122)
$ws.Range("B88").Value = 45676

# Rows 89-92 — four game results for 2025-01-19, formatted like row 80.
$ws.Range("B80:H80").Copy()
$ws.Range("B89:H89").PasteSpecial(-4122)
$ws.Range("B80:H80").Copy()
$ws.Range("B90:H90").PasteSpecial(-4122)
$ws.Range("B80:H80").Copy()
$ws.Range("B91:H91").PasteSpecial(-4122)
$ws.Range("B80:H80").Copy()
$ws.Range("B92:H92").PasteSpecial(-4122)

$ws.Rows.Item(89).RowHeight = 19.95
$ws.Rows.Item(90).RowHeight = 19.95
$ws.Rows.Item(91).RowHeight = 19.95
$ws.Rows.Item(92).RowHeight = 19.95

$ws.Range("B89").Value = "VSS - SIRIUS 43:84 (11:00, БНТУ)"
$ws.Range("B90").Value = "ISsoft - GOLDEN HILL 89:66 (12:30, БНТУ)"
$ws.Range("B91").Value = "Эра-Недвижимости плюс - ЛФК 72:42 (14:00, БНТУ)"
$ws.Range("B92").Value = "NORD - Mapogo males 64:90 (15:30, БНТУ)"

# ---------------------------------------------------------------------
# 3. Merge B:H on every new row, matching the rest of the sheet.
# ---------------------------------------------------------------------

$ws.Range("B84:H84").Merge()
$ws.Range("B85:H85").Merge()
$ws.Range("B86:H86").Merge()
$ws.Range("B87:H87").Merge()
$ws.Range("B88:H88").Merge()
$ws.Range("B89:H89").Merge()
$ws.Range("B90:H90").Merge()
$ws.Range("B91:H91").Merge()
$ws.Range("B92:H92").Merge()
